$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 07:19:11"
$ws.Range("A3").Value = "Total filas: 75"
$ws.Range("A49").Value = "07:19:11"
$ws.Range("B49").Value = "07:20"
$ws.Range("C49").Value = "16_SANTA ANA"
$ws.Range("D49").Value = 1
$ws.Range("A50").Value = "07:19:11"
$ws.Range("B50").Value = "07:20"
$ws.Range("D50").Value = 1
$ws.Range("A51").Value = "07:19:11"
$ws.Range("B51").Value = "07:21"
$ws.Range("C51").Value = "26_HERNANDEZ"
$ws.Range("D51").Value = 2
$ws.Range("A52").Value = "06:54:14"
$ws.Range("B52").Value = "07:23"
$ws.Range("C52").Value = "10_OLMOS"
$ws.Range("D52").Value = 29
$ws.Range("A53").Value = "05:55:25"
$ws.Range("B53").Value = "07:31"
$ws.Range("D53").Value = 96
$ws.Range("A54").Value = "05:55:25"
$ws.Range("B54").Value = "07:31"
$ws.Range("D54").Value = 96
$ws.Range("A55").Value = "07:19:11"
$ws.Range("D55").Value = 13
$ws.Range("A56").Value = "07:19:11"
$ws.Range("B56").Value = "07:32"
$ws.Range("C56").Value = "11_ETCHEVERRY"
$ws.Range("D56").Value = 13
$ws.Range("A57").Value = "07:19:11"
$ws.Range("B57").Value = "07:32"
$ws.Range("C57").Value = "16_SANTA ANA"
$ws.Range("D57").Value = 13
$ws.Range("A58").Value = "07:19:11"
$ws.Range("B58").Value = "07:35"
$ws.Range("C58").Value = "23_HERNANDEZ"
$ws.Range("D58").Value = 16
$ws.Range("B59").Value = "07:36"
$ws.Range("C59").Value = "27_EL RETIRO"
$ws.Range("D59").Value = 101
$ws.Range("A60").Value = "07:19:11"
$ws.Range("B60").Value = "07:37"
$ws.Range("C60").Value = "27_EL RETIRO"
$ws.Range("D60").Value = 18
$ws.Range("A61").Value = "07:19:11"
$ws.Range("B61").Value = "07:39"
$ws.Range("C61").Value = "10_OLMOS"
$ws.Range("D61").Value = 20
$ws.Range("A62").Value = "05:55:25"
$ws.Range("B62").Value = "07:47"
$ws.Range("C62").Value = "14_ABASTO"
$ws.Range("D62").Value = 112
$ws.Range("A63").Value = "07:19:11"
$ws.Range("B63").Value = "07:48"
$ws.Range("C63").Value = "14_ABASTO"
$ws.Range("D63").Value = 29
$ws.Range("B64").Value = "07:51"
$ws.Range("C64").Value = "215D_EL PATO"
$ws.Range("D64").Value = 57
$ws.Range("A65").Value = "07:19:11"
$ws.Range("B65").Value = "07:52"
$ws.Range("C65").Value = "215D_EL PATO"
$ws.Range("D65").Value = 33
$ws.Range("A66").Value = "07:19:11"
$ws.Range("B66").Value = "08:00"
$ws.Range("C66").Value = "23_HERNANDEZ"
$ws.Range("D66").Value = 41
$ws.Range("A67").Value = "06:26:08"
$ws.Range("B67").Value = "08:01"
$ws.Range("C67").Value = "23_HERNANDEZ"
$ws.Range("D67").Value = 95
$ws.Range("A68").Value = "07:19:11"
$ws.Range("B68").Value = "08:04"
$ws.Range("C68").Value = "11_ETCHEVERRY"
$ws.Range("D68").Value = 45
$ws.Range("B69").Value = "08:06"
$ws.Range("C69").Value = "23_HERNANDEZ"
$ws.Range("D69").Value = 72
$ws.Range("A70").Value = "07:19:11"
$ws.Range("B70").Value = "08:12"
$ws.Range("C70").Value = "15_ABASTO"
$ws.Range("D70").Value = 53
$ws.Range("A71").Value = "07:19:11"
$ws.Range("B71").Value = "08:21"
$ws.Range("C71").Value = "26_HERNANDEZ"
$ws.Range("D71").Value = 62
$ws.Range("E71").Value = "LP1912"
$ws.Range("A72").Value = "07:19:11"
$ws.Range("B72").Value = "08:23"
$ws.Range("C72").Value = "215B_EL PATO"
$ws.Range("D72").Value = 64
$ws.Range("E72").Value = "LP1912"
$ws.Range("A73").Value = "07:19:11"
$ws.Range("B73").Value = "08:23"
$ws.Range("C73").Value = "16_P MOR-SANTA ANA"
$ws.Range("D73").Value = 64
$ws.Range("E73").Value = "LP1912"
$ws.Range("A74").Value = "07:19:11"
$ws.Range("B74").Value = "08:27"
$ws.Range("C74").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D74").Value = 68
$ws.Range("E74").Value = "LP1912"
$ws.Range("A75").Value = "07:19:11"
$ws.Range("B75").Value = "08:42"
$ws.Range("C75").Value = "81_EL PELIGRO"
$ws.Range("D75").Value = 83
$ws.Range("E75").Value = "LP1912"
$ws.Range("A76").Value = "07:19:11"
$ws.Range("B76").Value = "08:44"
$ws.Range("C76").Value = "14_ABASTO"
$ws.Range("D76").Value = 85
$ws.Range("E76").Value = "LP1912"
$ws.Range("A77").Value = "07:19:11"
$ws.Range("B77").Value = "08:54"
$ws.Range("C77").Value = "17_ROMERO"
$ws.Range("D77").Value = 95
$ws.Range("E77").Value = "LP1912"
$ws.Range("A78").Value = "07:19:11"
$ws.Range("B78").Value = "09:02"
$ws.Range("C78").Value = "215A_EL PATO"
$ws.Range("D78").Value = 103
$ws.Range("E78").Value = "LP1912"
$ws.Range("A79").Value = "07:19:11"
$ws.Range("B79").Value = "09:11"
$ws.Range("C79").Value = "16_P MOR-SANTA ANA"
$ws.Range("D79").Value = 112
$ws.Range("E79").Value = "LP1912"
$ws.Range("A80").Value = "07:19:11"
$ws.Range("B80").Value = "09:17"
$ws.Range("C80").Value = "27_EL RETIRO"
$ws.Range("D80").Value = 118
$ws.Range("E80").Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 07:19:11"
$ws.Range("A3").Value = "Total filas: 14"
$ws.Range("A17").Value = "07:19:11"
$ws.Range("D17").Value = 33
$ws.Range("A18").Value = "07:19:11"
$ws.Range("D18").Value = 64
$ws.Range("A19").Value = "07:19:11"
$ws.Range("B19").Value = "09:02"
$ws.Range("C19").Value = "215A_EL PATO"
$ws.Range("D19").Value = 103
$ws.Range("E19").Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 07:19:11"
$ws.Range("A3").Value = "Total filas: 12"
$ws.Range("A13").Value = "07:19:11"
$ws.Range("D13").Value = 16
$ws.Range("A14").Value = "07:19:11"
$ws.Range("D14").Value = 48
$ws.Range("A16").Value = "07:19:11"
$ws.Range("B16").Value = "08:35"
$ws.Range("C16").Value = "215A_LA PLATA"
$ws.Range("D16").Value = 76
$ws.Range("E16").Value = "L6173"
$ws.Range("A17").Value = "07:19:11"
$ws.Range("B17").Value = "09:09"
$ws.Range("C17").Value = "215D_LA PLATA"
$ws.Range("D17").Value = 110
$ws.Range("E17").Value = "L6203"
